$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-89)
# from 45175 to 45177.
$ws.Range("C2:C89").Value2 = 45177

# Row 4 specific updates: a new species ("Klotsporig murkla") was added,
# which bumps the VU (K), Rödlistade (O), Hotade (P) and Alla arter (Q)
# counters and prepends the new name to the species list (R).
$ws.Range("K4").Value2 = 1
$ws.Range("O4").Value2 = 3
$ws.Range("P4").Value2 = 1
$ws.Range("Q4").Value2 = 9

$nl = "`r`n"
$ws.Range("R4").Value2 = "Klotsporig murkla${nl}Motaggsvamp${nl}Spillkråka${nl}Blåmossa${nl}Rödbrun jordstjärna${nl}Stubbspretmossa${nl}Svavelriska${nl}Thomsons trägnagare${nl}Blåsippa"
